# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped figures, mirroring the GitHub Actions refresh commit.
#
# The source cells are plain text (written by openpyxl as inline strings), so
# numeric-looking prices in column D must be forced to stay text instead of
# being auto-parsed into numbers by Excel's normal "typed input" handling.
# We do this the same way a human would in Excel: type a leading apostrophe
# to force text entry, then clear the resulting "quote prefix" number format
# so the cell's style is left exactly as it was (no stray formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values, keyed by row number.
$priceUpdates = @{
    2  = "27.636.73"
    3  = "1.616.58"
    5  = "209.14"
    8  = "23.08"
    12 = "1.846.93"
    13 = "1.624.75"
    16 = "64.54"
    17 = "27.650.12"
    18 = "226.52"
    19 = "7.64"
    25 = "154.20"
    26 = "6.87"
    28 = "15.39"
    29 = "0.992"
    32 = "3.36"
    34 = "1.389.26"
    35 = "1.58"
    36 = "0.996"
    39 = "0.553"
    40 = "0.840"
    44 = "65.39"
    46 = "1.756.45"
    48 = "87.50"
    51 = "7.56"
}

# New "Volume(1h)" (column E) values, keyed by row number.
$volumeUpdates = @{
    2  = "  -0.72%  "
    3  = "  -1.03%  "
    4  = "  -0.69%  "
    5  = "  -1.26%  "
    6  = "  -1.39%  "
    7  = "  -0.68%  "
    8  = "  -0.95%  "
    9  = "  -1.30%  "
    10 = "  -1.61%  "
    11 = "  -0.98%  "
    12 = "  -1.13%  "
    13 = "  -0.57%  "
    14 = "  -1.85%  "
    15 = "  -1.79%  "
    16 = "  -1.20%  "
    17 = "  -0.82%  "
    18 = "  -1.88%  "
    21 = "  -0.70%  "
    22 = "  -1.59%  "
    23 = "  -2.99%  "
    24 = "  -2.29%  "
    25 = "  -0.30%  "
    26 = "  -1.28%  "
    27 = "  -1.05%  "
    28 = "  -1.68%  "
    29 = "  -0.62%  "
    30 = "  -1.18%  "
    31 = "  -1.09%  "
    32 = "  -1.13%  "
    33 = "  -0.34%  "
    34 = "  -1.26%  "
    35 = "  +1.18%  "
    36 = "  -2.74%  "
    37 = "  -1.64%  "
    38 = "  -0.14%  "
    40 = "  -3.44%  "
    41 = "  -1.13%  "
    42 = "  -0.79%  "
    43 = "  -0.92%  "
    44 = "  -2.09%  "
    45 = "  -2.94%  "
    46 = "  -1.21%  "
    47 = "  -3.80%  "
    49 = "  +1.04%  "
    50 = "  -0.73%  "
    51 = "  +0.60%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $priceUpdates[$row]
    $cell.ClearFormats()
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
